$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 1).Value = "INDUSINDBK"
$ws1.Cells.Item(2, 2).Value = 801.85
$ws1.Cells.Item(2, 3).Value = 810
$ws1.Cells.Item(2, 4).Value = 789.1
$ws1.Cells.Item(2, 5).Value = 808.6
$ws1.Cells.Item(2, 6).Value = 2553303
$ws1.Cells.Item(2, 7).Value = 5249912
$ws1.Cells.Item(2, 8).Value = -0.5136484192496941
$ws1.Cells.Item(2, 9).Value = "INDUSINDBK"
$ws1.Cells.Item(3, 1).Value = "INDIGO"
$ws1.Cells.Item(3, 2).Value = 5893.5
$ws1.Cells.Item(3, 3).Value = 5907
$ws1.Cells.Item(3, 4).Value = 5793.5
$ws1.Cells.Item(3, 5).Value = 5884
$ws1.Cells.Item(3, 6).Value = 483125
$ws1.Cells.Item(3, 7).Value = 980678
$ws1.Cells.Item(3, 8).Value = -0.507356135245208
$ws1.Cells.Item(3, 9).Value = "INDIGO"
$ws1.Cells.Item(4, 1).Value = "BSE"
$ws1.Cells.Item(4, 2).Value = 2370
$ws1.Cells.Item(4, 3).Value = 2454
$ws1.Cells.Item(4, 4).Value = 2359.3
$ws1.Cells.Item(4, 5).Value = 2452
$ws1.Cells.Item(4, 6).Value = 4564283
$ws1.Cells.Item(4, 7).Value = 9401319
$ws1.Cells.Item(4, 8).Value = -0.5145061028138711
$ws1.Cells.Item(4, 9).Value = "BSE"
$ws1.Cells.Item(5, 1).Value = "POLYCAB"
$ws1.Cells.Item(5, 2).Value = 6790
$ws1.Cells.Item(5, 3).Value = 6903.5
$ws1.Cells.Item(5, 4).Value = 6746
$ws1.Cells.Item(5, 5).Value = 6896
$ws1.Cells.Item(5, 6).Value = 193378
$ws1.Cells.Item(5, 7).Value = 444662
$ws1.Cells.Item(5, 8).Value = -0.565112377491218
$ws1.Cells.Item(5, 9).Value = "POLYCAB"
$ws1.Cells.Item(6, 1).Value = "SBICARD"
$ws1.Cells.Item(6, 2).Value = 787.3
$ws1.Cells.Item(6, 3).Value = 799
$ws1.Cells.Item(6, 4).Value = 784
$ws1.Cells.Item(6, 5).Value = 799
$ws1.Cells.Item(6, 6).Value = 526708
$ws1.Cells.Item(6, 7).Value = 1239184
$ws1.Cells.Item(6, 8).Value = -0.5749557773502563
$ws1.Cells.Item(6, 9).Value = "SBICARD"
$ws1.Cells.Item(7, 1).Value = "MARICO"
$ws1.Cells.Item(7, 2).Value = 712.15
$ws1.Cells.Item(7, 3).Value = 719.9
$ws1.Cells.Item(7, 4).Value = 707.25
$ws1.Cells.Item(7, 5).Value = 718.1
$ws1.Cells.Item(7, 6).Value = 977299
$ws1.Cells.Item(7, 7).Value = 2396780
$ws1.Cells.Item(7, 8).Value = -0.5922450120578443
$ws1.Cells.Item(7, 9).Value = "MARICO"
$ws1.Cells.Item(8, 1).Value = "AMBER"
$ws1.Cells.Item(8, 2).Value = 7595
$ws1.Cells.Item(8, 3).Value = 7785
$ws1.Cells.Item(8, 4).Value = 7595
$ws1.Cells.Item(8, 5).Value = 7765
$ws1.Cells.Item(8, 6).Value = 185434
$ws1.Cells.Item(8, 7).Value = 397729
$ws1.Cells.Item(8, 8).Value = -0.5337679676362549
$ws1.Cells.Item(8, 9).Value = "AMBER"
$ws1.Cells.Item(9, 1).Value = "ANGELONE"
$ws1.Cells.Item(9, 2).Value = 2576
$ws1.Cells.Item(9, 3).Value = 2653
$ws1.Cells.Item(9, 4).Value = 2576
$ws1.Cells.Item(9, 5).Value = 2645.8
$ws1.Cells.Item(9, 6).Value = 543383
$ws1.Cells.Item(9, 7).Value = 1180598
$ws1.Cells.Item(9, 8).Value = -0.5397391830241962
$ws1.Cells.Item(9, 9).Value = "ANGELONE"
$ws1.Cells.Item(10, 1).Value = "IIFL"
$ws1.Cells.Item(10, 2).Value = 452
$ws1.Cells.Item(10, 3).Value = 455.4
$ws1.Cells.Item(10, 4).Value = 441.1
$ws1.Cells.Item(10, 5).Value = 453.35
$ws1.Cells.Item(10, 6).Value = 1488749
$ws1.Cells.Item(10, 7).Value = 3539104
$ws1.Cells.Item(10, 8).Value = -0.5793429636427752
$ws1.Cells.Item(10, 9).Value = "IIFL"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 1).Value = "TECHM"
$ws2.Cells.Item(2, 2).Value = 1455
$ws2.Cells.Item(2, 3).Value = 1492
$ws2.Cells.Item(2, 4).Value = 1450
$ws2.Cells.Item(2, 5).Value = 1488.1
$ws2.Cells.Item(2, 6).Value = 1587489
$ws2.Cells.Item(2, 7).Value = 1039973
$ws2.Cells.Item(2, 8).Value = 0.526471360314162
$ws2.Cells.Item(2, 9).Value = "TECHM"
$ws2.Cells.Item(3, 1).Value = "BAJFINANCE"
$ws2.Cells.Item(3, 2).Value = 871
$ws2.Cells.Item(3, 3).Value = 882
$ws2.Cells.Item(3, 4).Value = 867.25
$ws2.Cells.Item(3, 5).Value = 881.9
$ws2.Cells.Item(3, 6).Value = 5998002
$ws2.Cells.Item(3, 7).Value = 4212648
$ws2.Cells.Item(3, 8).Value = 0.4238080181396594
$ws2.Cells.Item(3, 9).Value = "BAJFINANCE"
$ws2.Cells.Item(4, 1).Value = "SUNPHARMA"
$ws2.Cells.Item(4, 2).Value = 1581
$ws2.Cells.Item(4, 3).Value = 1605.4
$ws2.Cells.Item(4, 4).Value = 1571
$ws2.Cells.Item(4, 5).Value = 1599.5
$ws2.Cells.Item(4, 6).Value = 3435674
$ws2.Cells.Item(4, 7).Value = 2376257
$ws2.Cells.Item(4, 8).Value = 0.4458343520923873
$ws2.Cells.Item(4, 9).Value = "SUNPHARMA"
$ws2.Cells.Item(5, 1).Value = "JIOFIN"
$ws2.Cells.Item(5, 2).Value = 323
$ws2.Cells.Item(5, 3).Value = 328
$ws2.Cells.Item(5, 4).Value = 316.8
$ws2.Cells.Item(5, 5).Value = 327.4
$ws2.Cells.Item(5, 6).Value = 16461908
$ws2.Cells.Item(5, 7).Value = 11615135
$ws2.Cells.Item(5, 8).Value = 0.4172808150744696
$ws2.Cells.Item(5, 9).Value = "JIOFIN"
$ws2.Cells.Item(6, 1).Value = "GRASIM"
$ws2.Cells.Item(6, 2).Value = 2764
$ws2.Cells.Item(6, 3).Value = 2772.7
$ws2.Cells.Item(6, 4).Value = 2718.8
$ws2.Cells.Item(6, 5).Value = 2741.9
$ws2.Cells.Item(6, 6).Value = 591589
$ws2.Cells.Item(6, 7).Value = 393906
$ws2.Cells.Item(6, 8).Value = 0.5018532340202992
$ws2.Cells.Item(6, 9).Value = "GRASIM"
$ws2.Cells.Item(7, 1).Value = "TORNTPHARM"
$ws2.Cells.Item(7, 2).Value = 3490.2
$ws2.Cells.Item(7, 3).Value = 3637
$ws2.Cells.Item(7, 4).Value = 3490.2
$ws2.Cells.Item(7, 5).Value = 3595
$ws2.Cells.Item(7, 6).Value = 402653
$ws2.Cells.Item(7, 7).Value = 257376
$ws2.Cells.Item(7, 8).Value = 0.564454339176924
$ws2.Cells.Item(7, 9).Value = "TORNTPHARM"
$ws2.Cells.Item(8, 1).Value = "PIDILITIND"
$ws2.Cells.Item(8, 2).Value = 3057
$ws2.Cells.Item(8, 3).Value = 3117.9
$ws2.Cells.Item(8, 4).Value = 3022.8
$ws2.Cells.Item(8, 5).Value = 3051
$ws2.Cells.Item(8, 6).Value = 1176910
$ws2.Cells.Item(8, 7).Value = 819116
$ws2.Cells.Item(8, 8).Value = 0.4368050434859043
$ws2.Cells.Item(8, 9).Value = "PIDILITIND"
$ws2.Cells.Item(9, 1).Value = "IRFC"
$ws2.Cells.Item(9, 2).Value = 127
$ws2.Cells.Item(9, 3).Value = 127.95
$ws2.Cells.Item(9, 4).Value = 124.1
$ws2.Cells.Item(9, 5).Value = 127.9
$ws2.Cells.Item(9, 6).Value = 9602587
$ws2.Cells.Item(9, 7).Value = 6331289
$ws2.Cells.Item(9, 8).Value = 0.5166875181341429
$ws2.Cells.Item(9, 9).Value = "IRFC"
$ws2.Cells.Item(10, 1).Value = "LUPIN"
$ws2.Cells.Item(10, 2).Value = 1855.9
$ws2.Cells.Item(10, 3).Value = 1952.7
$ws2.Cells.Item(10, 4).Value = 1855
$ws2.Cells.Item(10, 5).Value = 1937.3
$ws2.Cells.Item(10, 6).Value = 4166570
$ws2.Cells.Item(10, 7).Value = 2735877
$ws2.Cells.Item(10, 8).Value = 0.5229376174440591
$ws2.Cells.Item(10, 9).Value = "LUPIN"
$ws2.Cells.Item(11, 1).Value = "POLICYBZR"
$ws2.Cells.Item(11, 2).Value = 1730
$ws2.Cells.Item(11, 3).Value = 1790
$ws2.Cells.Item(11, 4).Value = 1725
$ws2.Cells.Item(11, 5).Value = 1786.3
$ws2.Cells.Item(11, 6).Value = 1072869
$ws2.Cells.Item(11, 7).Value = 710140
$ws2.Cells.Item(11, 8).Value = 0.5107851972850423
$ws2.Cells.Item(11, 9).Value = "POLICYBZR"
$ws2.Cells.Item(12, 1).Value = "IREDA"
$ws2.Cells.Item(12, 2).Value = 143.8
$ws2.Cells.Item(12, 3).Value = 145.45
$ws2.Cells.Item(12, 4).Value = 140.26
$ws2.Cells.Item(12, 5).Value = 145.45
$ws2.Cells.Item(12, 6).Value = 6734089
$ws2.Cells.Item(12, 7).Value = 4802241
$ws2.Cells.Item(12, 8).Value = 0.4022805186162044
$ws2.Cells.Item(12, 9).Value = "IREDA"
$ws2.Cells.Item(13, 1).Value = "ATGL"
$ws2.Cells.Item(13, 2).Value = 589
$ws2.Cells.Item(13, 3).Value = 594.75
$ws2.Cells.Item(13, 4).Value = 577.1
$ws2.Cells.Item(13, 5).Value = 591
$ws2.Cells.Item(13, 6).Value = 631605
$ws2.Cells.Item(13, 7).Value = 404411
$ws2.Cells.Item(13, 8).Value = 0.5617898622935578
$ws2.Cells.Item(13, 9).Value = "ATGL"
$ws2.Cells.Item(14, 1).Value = "SUZLON"
$ws2.Cells.Item(14, 2).Value = 64.81
$ws2.Cells.Item(14, 3).Value = 65.26
$ws2.Cells.Item(14, 4).Value = 62.21
$ws2.Cells.Item(14, 5).Value = 64.2
$ws2.Cells.Item(14, 6).Value = 75542321
$ws2.Cells.Item(14, 7).Value = 47280146
$ws2.Cells.Item(14, 8).Value = 0.597759892704223
$ws2.Cells.Item(14, 9).Value = "SUZLON"
$ws2.Cells.Item(15, 1).Value = "SRF"
$ws2.Cells.Item(15, 2).Value = 2965
$ws2.Cells.Item(15, 3).Value = 2965
$ws2.Cells.Item(15, 4).Value = 2875
$ws2.Cells.Item(15, 5).Value = 2913
$ws2.Cells.Item(15, 6).Value = 479056
$ws2.Cells.Item(15, 7).Value = 321329
$ws2.Cells.Item(15, 8).Value = 0.4908582792091595
$ws2.Cells.Item(15, 9).Value = "SRF"
$ws2.Cells.Item(16, 1).Value = "MANAPPURAM"
$ws2.Cells.Item(16, 2).Value = 255.6
$ws2.Cells.Item(16, 3).Value = 263.15
$ws2.Cells.Item(16, 4).Value = 253.7
$ws2.Cells.Item(16, 5).Value = 261.9
$ws2.Cells.Item(16, 6).Value = 2066794
$ws2.Cells.Item(16, 7).Value = 1423890
$ws2.Cells.Item(16, 8).Value = 0.4515124061549698
$ws2.Cells.Item(16, 9).Value = "MANAPPURAM"

